$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 119; everything from row 119 down shifts to row 120+.
$ws.Rows("119:119").Insert()

# Populate the newly inserted row 119 with the new weekly price-report entry.
$ws.Range("A119").Value = 7
$ws.Range("B119").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C119").Value = "Ñuble"
$ws.Range("D119").Value = 44468
$ws.Range("E119").Value = 16
$ws.Range("F119").Value = 100112009
$ws.Range("G119").Value = "Acelga"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 100
$ws.Range("K119").Value = 400
$ws.Range("L119").Value = 450
$ws.Range("M119").Value = 425
$ws.Range("N119").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O119").Value = "Provincia de Diguillín"
$ws.Range("P119").Value = 425
$ws.Range("Q119").Value = 1
$ws.Range("R119").Value = "Hortaliza"
